# Clash-of-clans workbook: rename / split the single "Sheet3" data dump into
# three tabs (Troops, Defense, Heros), refresh troop stats, add a new troop
# (Ice Golem) and populate brand-new Defense + Heros sheets.
#
# NOTE: the order in which brand-new text values are first written below is
# deliberate (it mirrors the order new entries were appended to the shared
# string table in the original edit) - don't reorder the blocks casually.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing sheet to "Troops" and add the two new sheets
#    right after it, in tab order Troops -> Defense -> Heros.
# ---------------------------------------------------------------------
$troops = $wb.Worksheets.Item(1)
$troops.Name = "Troops"

$defense = $wb.Worksheets.Add($null, $troops)
$defense.Name = "Defense"

$heros = $wb.Worksheets.Add($null, $defense)
$heros.Name = "Heros"

# ---------------------------------------------------------------------
# 2. Troops sheet header + the 7 troops whose names/columns don't change.
# ---------------------------------------------------------------------
$troops.Range("A1").Value = "troop"
$troops.Range("B1").Value = "housing"
$troops.Range("C1").Value = "damage_per_second"
$troops.Range("D1").Value = "hitpoints"
$troops.Range("E1").Value = "cost"

$troopRows = @(
    @(2,  "Barbarian",     1,  34,  205,  300),
    @(3,  "Archer",        1,  28,  52,   600),
    @(4,  "Giant",         5,  64,  1660, 4500),
    @(5,  "Goblin",        1,  52,  101,  200),
    @(6,  "Wall Breaker",  2,  100, 70,   2750),
    @(7,  "Balloon",       5,  236, 840,  5500),
    @(8,  "Wizard",        4,  230, 230,  5500)
)
foreach ($r in $troopRows) {
    $row = $r[0]
    $troops.Cells.Item($row, 1).Value = $r[1]
    $troops.Cells.Item($row, 2).Value = $r[2]
    $troops.Cells.Item($row, 3).Value = $r[3]
    $troops.Cells.Item($row, 4).Value = $r[4]
    $troops.Cells.Item($row, 5).Value = $r[5]
}

# New 21st troop row - first brand-new string in the workbook.
$troops.Cells.Item(22, 1).Value = "Ice Golem"
$troops.Cells.Item(22, 2).Value = 15
$troops.Cells.Item(22, 3).Value = 36
$troops.Cells.Item(22, 4).Value = 3200
$troops.Cells.Item(22, 5).Value = 28000

# ---------------------------------------------------------------------
# 3. Heros sheet data rows, then its header (range header text is
#    deferred - see step 6).
# ---------------------------------------------------------------------
$herosRows = @(
    @(2, "Archer Queen",   658, 2575, 5),
    @(3, "Barbarian King", 410, 7303, 0),
    @(4, "Grand Warden",   168, 1833, 7)
)
foreach ($r in $herosRows) {
    $row = $r[0]
    $heros.Cells.Item($row, 1).Value = $r[1]
    $heros.Cells.Item($row, 2).Value = $r[2]
    $heros.Cells.Item($row, 3).Value = $r[3]
    # Column D is numeric (not a new shared string), safe to write now.
    $heros.Cells.Item($row, 4).Value = $r[4]
}
$heros.Range("A1").Value = "hero"
$heros.Range("B1").Value = "damage_per_second"
$heros.Range("C1").Value = "hitpoints"

# ---------------------------------------------------------------------
# 4. Defense sheet header + first row (Cannon).
# ---------------------------------------------------------------------
$defense.Range("A1").Value = "building"
$defense.Range("B1").Value = "damage_per_second"
$defense.Range("C1").Value = "hitpoints"

$defense.Cells.Item(2, 1).Value = "Cannon"
$defense.Cells.Item(2, 2).Value = 125
$defense.Cells.Item(2, 3).Value = 1620
$defense.Cells.Item(2, 4).Value = 9

# ---------------------------------------------------------------------
# 5. Troops sheet: re-write rows 9-21 with proper-cased names and
#    refreshed stats.
# ---------------------------------------------------------------------
$troopRows2 = @(
    @(9,  "Healer",         14, 80,  1500, 15000),
    @(10, "Dragon",         20, 310, 3600, 30000),
    @(11, "Pekka",          25, 610, 6300, 39000),
    @(12, "Baby Dragon",    10, 125, 1700, 15000),
    @(13, "Miner",          6,  120, 870,  6400),
    @(14, "Minion",         2,  62,  96,   1300),
    @(15, "Hog Rider",      5,  148, 810,  14000),
    @(16, "Valkyrie",       8,  178, 1300, 25000),
    @(17, "Golem",          30, 66,  7200, 82500),
    @(18, "Witch",          12, 160, 440,  32500),
    @(19, "Lava Hound",     30, 18,  7600, 63000),
    @(20, "Bowler",         6,  90,  390,  17000),
    @(21, "Electro Dragon", 30, 300, 4200, 44000)
)
foreach ($r in $troopRows2) {
    $row = $r[0]
    $troops.Cells.Item($row, 1).Value = $r[1]
    $troops.Cells.Item($row, 2).Value = $r[2]
    $troops.Cells.Item($row, 3).Value = $r[3]
    $troops.Cells.Item($row, 4).Value = $r[4]
    $troops.Cells.Item($row, 5).Value = $r[5]
}

# Rows 17 & 19 used to hold a formula for hitpoints; they are now plain
# values, so make sure no stale formula remains.
$troops.Cells.Item(17, 4).Value = 7200
$troops.Cells.Item(19, 4).Value = 7600

# ---------------------------------------------------------------------
# 6. Defense sheet row 3 (Archer Tower), then the "Range" header text
#    (first used here, then re-used verbatim on the Heros sheet).
# ---------------------------------------------------------------------
$defense.Cells.Item(3, 1).Value = "Archer Tower"
$defense.Cells.Item(3, 2).Value = 116
$defense.Cells.Item(3, 3).Value = 1330
$defense.Cells.Item(3, 4).Value = 10

$defense.Range("D1").Value = "Range"
$heros.Range("D1").Value = "Range"

# ---------------------------------------------------------------------
# 7. Remaining Defense sheet rows (4-12).
# ---------------------------------------------------------------------
$defenseRows = @(
    @(4,  "Mortar",          35,  900,  11),
    @(5,  "Air Defense",     400, 1400, 10),
    @(6,  "Wizard Tower",    70,  2240, 7),
    @(7,  "Tesla",           120, 950,  7),
    @(8,  "Bomber Tower",    60,  1400, 6),
    @(9,  "X-Bow",           155, 3500, 14),
    @(10, "Inferno Tower",   81,  3000, 10),
    @(11, "Eagle Artillery", 120, 4800, 50),
    @(12, "Giga Tesla",      200, 7500, 10)
)
foreach ($r in $defenseRows) {
    $row = $r[0]
    $defense.Cells.Item($row, 1).Value = $r[1]
    $defense.Cells.Item($row, 2).Value = $r[2]
    $defense.Cells.Item($row, 3).Value = $r[3]
    $defense.Cells.Item($row, 4).Value = $r[4]
}

# ---------------------------------------------------------------------
# 8. Column widths / bold header formatting that came along with the
#    split. (ColumnWidth is specified in "characters"; the host rounds
#    the stored worksheet width to the nearest 1/6 character, so the
#    inputs below are pre-compensated to land as close as possible on
#    the target widths of 14 / 19.5703125 / 13.85546875.)
# ---------------------------------------------------------------------
$troops.Columns.Item(1).ColumnWidth = 13.166666666666666
$troops.Range("A1:E1").Font.Bold = $true

$defense.Columns.Item(2).ColumnWidth = 18.666666666666668
$defense.Range("A1:D1").Font.Bold = $true

$heros.Columns.Item(1).ColumnWidth = 13.0
$heros.Columns.Item(2).ColumnWidth = 18.666666666666668
$heros.Range("A1:D1").Font.Bold = $true

# ---------------------------------------------------------------------
# 9. Selection / active-tab bookkeeping to match the final saved state.
# ---------------------------------------------------------------------
$defense.Range("C13").Select()
$heros.Range("C5").Select()
$heros.Activate()
